$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计").
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "2022-Q1"

# Reuse the header/index-column formatting already present in the workbook
# (copy from the "2021-Q4" sheet, which uses the same layout/style).
$afterSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$afterSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$afterSheet.Range("A2").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)
$afterSheet.Range("A2").Copy()
$newSheet.Range("A4").PasteSpecial(-4122)
$afterSheet.Range("A2").Copy()
$newSheet.Range("A5").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Text-format the data columns (B:G) so fund codes / percentages keep
# their original textual form instead of being coerced to numbers.
$newSheet.Range("B2:G5").NumberFormat = "@"

# Row 2 - 501305
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "501305"
$newSheet.Range("C2").Value = "汇添富中证港股通高股息投资指数（LOF）A"
$newSheet.Range("D2").Value = "1.59"
$newSheet.Range("E2").Value = "93.08"
$newSheet.Range("F2").Value = "5.40"
$newSheet.Range("G2").Value = "0.0859"
$newSheet.Range("H2").Value = 2

# Row 3 - 501306
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "501306"
$newSheet.Range("C3").Value = "汇添富中证港股通高股息投资指数（LOF）C"
$newSheet.Range("D3").Value = "0.21"
$newSheet.Range("E3").Value = "93.08"
$newSheet.Range("F3").Value = "5.40"
$newSheet.Range("G3").Value = "0.0113"
$newSheet.Range("H3").Value = 2

# Row 4 - 501307
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "501307"
$newSheet.Range("C4").Value = "银河中证沪港深高股息指数（LOF）A"
$newSheet.Range("D4").Value = "0.19"
$newSheet.Range("E4").Value = "91.35"
$newSheet.Range("F4").Value = "2.11"
$newSheet.Range("G4").Value = "0.0040"
$newSheet.Range("H4").Value = 2

# Row 5 - 501308
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "501308"
$newSheet.Range("C5").Value = "银河中证沪港深高股息指数（LOF）C"
$newSheet.Range("D5").Value = "0.01"
$newSheet.Range("E5").Value = "91.35"
$newSheet.Range("F5").Value = "2.11"
$newSheet.Range("G5").Value = "0.0002"
$newSheet.Range("H5").Value = 2

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing rows down and renumbering the index column.
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Extend the index-column style (already on A2) down through A5.
$zj.Range("A2").Copy()
$zj.Range("A3").PasteSpecial(-4122)
$zj.Range("A2").Copy()
$zj.Range("A4").PasteSpecial(-4122)
$zj.Range("A2").Copy()
$zj.Range("A5").PasteSpecial(-4122)

# Former row 4 (2021-Q1) -> row 5
$zj.Range("A5").Value = 3
$zj.Range("B5").Value = "2021-Q1"
$zj.Range("C5").Value = 4
$zj.Range("D5").Value = 0.74

# Former row 3 (2021-Q3) -> row 4
$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2021-Q3"
$zj.Range("C4").Value = 2
$zj.Range("D4").Value = 0.03

# Former row 2 (2021-Q4) -> row 3
$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2021-Q4"
$zj.Range("C3").Value = 8
$zj.Range("D3").Value = 0.18

# New row 2 (2022-Q1)
$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 4
$zj.Range("D2").Value = 0.1

# Restore the originally active sheet/tab ("2021-Q1"), since adding the new
# worksheet made it the active one as a side effect.
$wb.Worksheets.Item("2021-Q1").Activate()
